$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.437.79"
$ws.Range("E2").Value = "  -6.81%  "

$ws.Range("D3").Value = "3.739.18"
$ws.Range("E3").Value = "  -5.93%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.74"
$ws.Range("E5").Value = "  -5.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.74"
$ws.Range("E6").Value = "  +5.59%  "

$ws.Range("D7").Value = "3.730.60"
$ws.Range("E7").Value = "  -6.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.636"
$ws.Range("E8").Value = "  -6.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.721"
$ws.Range("E10").Value = "  -5.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.167"
$ws.Range("E11").Value = "  -10.54%  "

$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000301"
$ws.Range("E13").Value = "  -10.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.85"
$ws.Range("E14").Value = "  -2.84%  "

$ws.Range("D15").Value = "4.339.58"
$ws.Range("E15").Value = "  -5.80%  "

$ws.Range("D16").Value = "3.763.78"
$ws.Range("E16").Value = "  -5.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.54"
$ws.Range("E17").Value = "  -4.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.14"
$ws.Range("E18").Value = "  -6.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.15"
$ws.Range("E19").Value = "  -6.86%  "

$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").Value = "68.327.83"
$ws.Range("E21").Value = "  -6.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "413.11"
$ws.Range("E22").Value = "  -5.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.59"
$ws.Range("E23").Value = "  -5.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.19"
$ws.Range("E24").Value = "  -7.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  -7.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.00"
$ws.Range("E26").Value = "  -8.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.86"
$ws.Range("E28").Value = "  -5.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.99"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.68"
$ws.Range("E30").Value = "  -8.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.11"
$ws.Range("E31").Value = "  +4.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.20"
$ws.Range("E32").Value = "  -8.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.83"
$ws.Range("E33").Value = "  -6.73%  "

$ws.Range("E34").Value = "  -7.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.45"
$ws.Range("E35").Value = "  -7.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "66.13"
$ws.Range("E36").Value = "  -6.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "619.35"
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").Value = "0.0₃0927"
$ws.Range("E38").Value = "  -12.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.405"
$ws.Range("E39").Value = "  -6.01%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.23"
$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.138"
$ws.Range("E43").Value = "  -5.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.11"
$ws.Range("E44").Value = "  -8.63%  "

$ws.Range("E45").Value = "  -7.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.65"
$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.53"
$ws.Range("E47").Value = "  -9.92%  "

$ws.Range("E48").Value = "  -14.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("E49").Value = "  -8.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -8.27%  "

$ws.Range("D51").Value = "2.748.29"
$ws.Range("E51").Value = "  -2.23%  "
